{"js": "// Replace the multiplication-answer strings in the table cells with the\n// new values from the commit. Each old value is unique in the document,\n// so a body-wide search/replace keyed on the exact old text is safe.\nconst replacements = [\n  [\"72\u00d738=2736\", \"74\u00d724=1776\"],\n  [\"11\u00d720=220\", \"46\u00d737=1702\"],\n  [\"13\u00d799=1287\", \"67\u00d717=1139\"],\n  [\"41\u00d798=4018\", \"24\u00d746=1104\"],\n  [\"36\u00d720=720\", \"46\u00d715=690\"],\n  [\"49\u00d720=980\", \"30\u00d729=870\"],\n  [\"91\u00d765=5915\", \"39\u00d773=2847\"],\n  [\"85\u00d784=7140\", \"42\u00d762=2604\"],\n  [\"45\u00d747=2115\", \"21\u00d798=2058\"],\n  [\"67\u00d735=2345\", \"26\u00d773=1898\"],\n  [\"42\u00d744=1848\", \"94\u00d780=7520\"],\n  [\"71\u00d775=5325\", \"21\u00d770=1470\"],\n  [\"59\u00d726=1534\", \"50\u00d785=4250\"],\n  [\"37\u00d788=3256\", \"59\u00d780=4720\"],\n  [\"49\u00d711=539\", \"87\u00d765=5655\"],\n  [\"93\u00d715=1395\", \"98\u00d781=7938\"],\n  [\"20\u00d733=660\", \"90\u00d798=8820\"],\n  [\"84\u00d733=2772\", \"73\u00d712=876\"],\n  [\"43\u00d749=2107\", \"20\u00d729=580\"],\n  [\"91\u00d743=3913\", \"99\u00d765=6435\"],\n  [\"41\u00d766=2706\", \"61\u00d721=1281\"],\n  [\"32\u00d754=1728\", \"33\u00d764=2112\"],\n  [\"66\u00d771=4686\", \"97\u00d755=5335\"],\n  [\"41\u00d724=984\", \"39\u00d793=3627\"],\n  [\"74\u00d726=1924\", \"31\u00d716=496\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-answer strings in the table cells with the\n# new values from the commit. Each old value is unique in the document,\n# so Find/Replace (Replace:=wdReplaceAll, value 2) keyed on the exact old\n# text is safe and touches only the intended run.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"72\u00d738=2736\", \"74\u00d724=1776\"),\n  @(\"11\u00d720=220\", \"46\u00d737=1702\"),\n  @(\"13\u00d799=1287\", \"67\u00d717=1139\"),\n  @(\"41\u00d798=4018\", \"24\u00d746=1104\"),\n  @(\"36\u00d720=720\", \"46\u00d715=690\"),\n  @(\"49\u00d720=980\", \"30\u00d729=870\"),\n  @(\"91\u00d765=5915\", \"39\u00d773=2847\"),\n  @(\"85\u00d784=7140\", \"42\u00d762=2604\"),\n  @(\"45\u00d747=2115\", \"21\u00d798=2058\"),\n  @(\"67\u00d735=2345\", \"26\u00d773=1898\"),\n  @(\"42\u00d744=1848\", \"94\u00d780=7520\"),\n  @(\"71\u00d775=5325\", \"21\u00d770=1470\"),\n  @(\"59\u00d726=1534\", \"50\u00d785=4250\"),\n  @(\"37\u00d788=3256\", \"59\u00d780=4720\"),\n  @(\"49\u00d711=539\", \"87\u00d765=5655\"),\n  @(\"93\u00d715=1395\", \"98\u00d781=7938\"),\n  @(\"20\u00d733=660\", \"90\u00d798=8820\"),\n  @(\"84\u00d733=2772\", \"73\u00d712=876\"),\n  @(\"43\u00d749=2107\", \"20\u00d729=580\"),\n  @(\"91\u00d743=3913\", \"99\u00d765=6435\"),\n  @(\"41\u00d766=2706\", \"61\u00d721=1281\"),\n  @(\"32\u00d754=1728\", \"33\u00d764=2112\"),\n  @(\"66\u00d771=4686\", \"97\u00d755=5335\"),\n  @(\"41\u00d724=984\", \"39\u00d793=3627\"),\n  @(\"74\u00d726=1924\", \"31\u00d716=496\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n"}
